$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (outside the used A1:E51 range) used to stage numeric-looking
# text so PasteSpecial(values) preserves it as text instead of Excel auto-
# converting the typed string to a real number via Range.Value.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$ws.Range("D2").Value = '61.486.23'
$ws.Range("E2").Value = '  -2.06%  '
$ws.Range("D3").Value = '3.371.64'
$ws.Range("E3").Value = '  -2.74%  '
$ws.Range("E4").Value = '  -0.01%  '
$scratch.Value = '403.60'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -2.53%  '
$scratch.Value = '132.20'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +6.14%  '
$scratch.Value = '0.585'
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -0.94%  '
$ws.Range("E8").Value = '  -0.01%  '
$scratch.Value = '0.665'
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  -2.62%  '
$scratch.Value = '0.118'
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -7.43%  '
$scratch.Value = '42.03'
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("D13").Value = '3.901.91'
$ws.Range("E13").Value = '  -2.89%  '
$scratch.Value = '8.36'
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  -2.59%  '
$scratch.Value = '19.70'
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  -1.25%  '
$ws.Range("D16").Value = '3.381.67'
$ws.Range("E16").Value = '  -2.66%  '
$ws.Range("D17").Value = '61.480.77'
$ws.Range("E17").Value = '  -2.06%  '
$ws.Range("E18").Value = '  -1.97%  '
$scratch.Value = '10.87'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +0.63%  '
$scratch.Value = '0.0000127'
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -7.21%  '
$scratch.Value = '3.18'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -4.35%  '
$scratch.Value = '84.67'
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +3.13%  '
$scratch.Value = '314.01'
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -0.75%  '
$scratch.Value = '12.68'
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  -2.03%  '
$scratch.Value = '3.09'
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  -2.63%  '
$scratch.Value = '4.77'
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  +10.48%  '
$scratch.Value = '29.31'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -4.92%  '
$scratch.Value = '8.18'
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +4.32%  '
$scratch.Value = '7.66'
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -1.78%  '
$scratch.Value = '2.67'
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +0.98%  '
$scratch.Value = '0.171'
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -2.16%  '
$ws.Range("E32").Value = '  -2.31%  '
$ws.Range("B33").Value = 'Dai'
$ws.Range("C33").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$scratch.Value = '0.999'
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -0.24%  '
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$scratch.Value = '41.60'
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -1.16%  '
$ws.Range("B35").Value = 'Cosmos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$scratch.Value = '11.29'
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -2.19%  '
$scratch.Value = '0.0477'
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -3.20%  '
$scratch.Value = '51.64'
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -1.03%  '
$scratch.Value = '0.999'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +0.12%  '
$scratch.Value = '3.39'
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -3.15%  '
$scratch.Value = '2.92'
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -3.79%  '
$scratch.Value = '138.51'
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +1.77%  '
$ws.Range("E42").Value = '  -1.48%  '
$scratch.Value = '0.123'
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -1.18%  '
$scratch.Value = '0.291'
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +2.44%  '
$scratch.Value = '3.95'
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +1.39%  '
$scratch.Value = '16.59'
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -1.83%  '
$ws.Range("E47").Value = '  -1.49%  '
$scratch.Value = '21.24'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").Value = '2.114.48'
$ws.Range("E49").Value = '  -4.15%  '
$ws.Range("E50").Value = '  -7.36%  '
$scratch.Value = '1.86'
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +0.78%  '

# Remove the scratch column entirely so it does not widen the sheets used range
$ws.Columns.Item(26).Delete()
$excel.CutCopyMode = 0
